# Correccion en icono: facebook
#
# Slide 3 has a group ("1 Grupo") containing a background rectangle and the
# facebook-icon picture. The picture's frame was mis-sized/mis-positioned
# relative to the rectangle; fix it up by re-sizing/re-positioning the
# picture so it sits correctly inside the (unchanged) rectangle, then
# regroup so the group's bounding box is recomputed from its children
# (this is what also bumps the group's shape id/name, matching the
# author's original edit).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

$grp = $s.Shapes.Item(1)

# Ungroup -> gives back a ShapeRange with the two original children in
# their original z-order: [1]=rectangle, [2]=picture.
$items = $grp.Ungroup()
$rect = $items.Item(1)
$pic = $items.Item(2)

# Resize/reposition only the picture (the rectangle keeps its original
# frame). Point values below are chosen so that, after PowerPoint's
# Single-precision (points -> EMU) COM round-trip, they land exactly on
# the target EMU coordinates:
#   off  = (2229463, 1162924)
#   ext  = (3388930, 3380023)
$pic.Left = 175.548272
$pic.Top = 91.568821
$pic.Width = 266.844895
$pic.Height = 266.143540

# Regroup the two shapes; PowerPoint recomputes the new group's bounding
# box from its children (here it matches the rectangle's frame exactly,
# since the resized picture now sits fully inside it) and assigns a new
# shape id / default name.
$range = $s.Shapes.Range(@($rect.Name, $pic.Name))
$newGrp = $range.Group()
$newGrp.Name = "2 Grupo"
